# "updated methods section, tabs & figs"
# Rebuild the "dataset" sheet (sheet2) with the new results table:
#   - 3 data-source columns (CBS / CANPATH / APL) with an "N (%)" sub-header
#   - Age Groups breakdown (18-35 / 36-64 / 65+) with combined "n (%)" text
#   - Gender / Provinces / Data Collecting Period row labels added
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset")

# Start from a clean slate so stale shared strings / cells don't linger.
$ws.Range("A1:D11").ClearContents()

# --- Row 3-6: Age Groups block (write label strings first so the shared
#     string table fills in the same order as the source workbook) -------
$ws.Range("A3").Value = "Age Groups"
$ws.Range("A4").Value = "18-35"
$ws.Range("A5").Value = "36-64"
$ws.Range("A6").Value = "65+"

# --- Row 1: data-source headers -----------------------------------------
$ws.Range("B1").Value = "CBS"
$ws.Range("C1").Value = "CANPATH"
$ws.Range("D1").Value = "APL"

# --- Rows 7-11: additional demographic / metadata row labels -------------
$ws.Range("A7").Value = "Gender"
$ws.Range("A8").Value = "Male"
$ws.Range("A9").Value = "Female"
$ws.Range("A10").Value = "Provinces"
$ws.Range("A11").Value = "Data Collecting Period"

# --- Row 2: "N (%)" sub-header under each data source ---------------------
$ws.Range("B2").Value = "N (%)"
$ws.Range("C2").Value = "N (%)"
$ws.Range("D2").Value = "N (%)"

# --- Rows 4-6: combined N (%) values for each age group -------------------
$ws.Range("B4").Value = "200 (33.33%)"
$ws.Range("B5").Value = "201 (33.33%)"
$ws.Range("B6").Value = "202 (33.33%)"

# --- Column widths: widen col A for the longer labels, size the two new
#     data columns -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.8333333333333
$ws.Columns.Item(3).ColumnWidth = 11.3333333333333
$ws.Columns.Item(4).ColumnWidth = 10.5

# --- Header row alignment tweak (B1:D1) -----------------------------------
$ws.Range("B1:D1").WrapText = $false

# --- Selection moves to B7 ------------------------------------------------
$ws.Activate()
$ws.Range("B7").Select()
